$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: B4 changes from text to a real number
$ws.Range("B4").Value = 2304221520006

# Row 5 (new)
$ws.Range("A5").Value = "abhi"
$ws.Range("B5").Value = 230
$ws.Range("C5").Value = "btech"
$ws.Range("D5").Value = "ai"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2"
$ws.Range("E5").Style = "Normal"

# Row 6 (new)
$ws.Range("A6").Value = "abhi"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "230"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "b.tech"
$ws.Range("D6").Value = "ai"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2"
$ws.Range("E6").Style = "Normal"
